# Update countries & provincias Spain
#
# The source COVID-19 table ("Pais" sheet) is kept sorted by column B
# ("Casos totales") descending. This edit refreshes the timestamp plus a
# handful of countries' figures; because some of the new totals change
# their rank, several rows swap places with their neighbour (the row
# keeps the OTHER country's untouched numbers while the updated country's
# row gets the fresh figures) and a couple of groups of tied rows simply
# rotate their country labels with no numeric change at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Page timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 01:05"

# Estados Unidos (row 4) - refreshed totals, stays #1
$ws.Range("B4").Value = 1569659
$ws.Range("C4").Value = 19365
$ws.Range("D4").Value = 362984
$ws.Range("E4").Value = 1113202
$ws.Range("G4").Value = 1492
$ws.Range("H4").Value = 93473

# Sudan overtakes Senegal (rows 78/79 swap; Sudan gets new figures,
# Senegal's figures are carried over unchanged)
$ws.Range("A78").Value = "Sudan"
$ws.Range("B78").Value = 2728
$ws.Range("C78").Value = 137
$ws.Range("D78").Value = 286
$ws.Range("E78").Value = 2331
$ws.Range("G78").Value = 6
$ws.Range("H78").Value = 111

$ws.Range("A79").Value = "Senegal"
$ws.Range("B79").Value = 2617
$ws.Range("C79").Value = 73
$ws.Range("D79").Value = 1133
$ws.Range("E79").Value = 1454
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = 30

# Guyana overtakes Bermudas (rows 164/165 swap; Guyana gets new figures,
# Bermudas' figures are carried over unchanged)
$ws.Range("A164").Value = "Guyana"
$ws.Range("C164").Value = 1
$ws.Range("D164").Value = 46
$ws.Range("E164").Value = 69
$ws.Range("H164").Value = 10

$ws.Range("A165").Value = "Bermudas"
$ws.Range("B165").Value = 125
$ws.Range("D165").Value = 77
$ws.Range("E165").Value = 39
$ws.Range("H165").Value = 9

# Montserrat / Groenlandia / Seychelles (rows 209-211) rotate labels;
# only Montserrat's own figures move with it, the other two are ties
# with identical totals so only the label moves
$ws.Range("A209").Value = "Montserrat"
$ws.Range("D209").Value = 10
$ws.Range("H209").Value = 1

$ws.Range("A210").Value = "Groenlandia"

$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Bonaire, San Eustaquio y Saba / San Bartolome / Sahara Occidental
# (rows 214-216) are a 3-way tie (identical figures) that simply rotate
# labels
$ws.Range("A214").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Sahara Occidental"
